$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 23.73148533333334
$ws.Range("H2").Value = 71.194456
$ws.Range("I2").Value = 0.8653076146801144
$ws.Range("J2").Value = 0.8653076146801145
$ws.Range("M2").Value = 19.163974
$ws.Range("N2").Value = 57.491922
$ws.Range("O2").Value = 0.6845732287637933
$ws.Range("P2").Value = 0.6845732287637933
$ws.Range("Q2").Value = 454.7895679093814
$ws.Range("R2").Value = 4093.106111184432
$ws.Range("S2").Value = 0.5923664276554623
$ws.Range("T2").Value = 0.5923664276554623
$ws.Range("G3").Value = 23.73148533333334
$ws.Range("H3").Value = 71.194456
$ws.Range("I3").Value = 0.8653076146801144
$ws.Range("J3").Value = 0.8653076146801145
$ws.Range("O3").Value = 0.02733363438148322
$ws.Range("P3").Value = 0.02733363438148323
$ws.Range("Q3").Value = 18.15883421587467
$ws.Range("R3").Value = 163.429507942872
$ws.Range("S3").Value = 0.02365200196717961
$ws.Range("T3").Value = 0.02365200196717962
$ws.Range("G4").Value = 23.73148533333334
$ws.Range("H4").Value = 71.194456
$ws.Range("I4").Value = 0.8653076146801144
$ws.Range("J4").Value = 0.8653076146801145
$ws.Range("M4").Value = 7.880893333333333
$ws.Range("N4").Value = 23.64268
$ws.Range("O4").Value = 0.281520346184098
$ws.Range("P4").Value = 0.281520346184098
$ws.Range("Q4").Value = 187.0253045535644
$ws.Range("R4").Value = 1683.22774098208
$ws.Range("S4").Value = 0.2436016992404819
$ws.Range("T4").Value = 0.2436016992404819
$ws.Range("G5").Value = 23.73148533333334
$ws.Range("H5").Value = 71.194456
$ws.Range("I5").Value = 0.8653076146801144
$ws.Range("J5").Value = 0.8653076146801145
$ws.Range("M5").Value = 0.183999
$ws.Range("N5").Value = 0.551997
$ws.Range("O5").Value = 0.006572790670625477
$ws.Range("P5").Value = 0.006572790670625476
$ws.Range("Q5").Value = 4.366569569848
$ws.Range("R5").Value = 39.299126128632
$ws.Range("S5").Value = 0.005687485816990641
$ws.Range("T5").Value = 0.005687485816990641
$ws.Range("I6").Value = 0.09010639372350319
$ws.Range("J6").Value = 0.09010639372350321
$ws.Range("M6").Value = 19.163974
$ws.Range("N6").Value = 57.491922
$ws.Range("O6").Value = 0.6845732287637933
$ws.Range("P6").Value = 0.6845732287637933
$ws.Range("Q6").Value = 47.358242516488
$ws.Range("R6").Value = 426.2241826483921
$ws.Range("S6").Value = 0.06168442488356018
$ws.Range("T6").Value = 0.06168442488356019
$ws.Range("I7").Value = 0.09010639372350319
$ws.Range("J7").Value = 0.09010639372350321
$ws.Range("O7").Value = 0.02733363438148322
$ws.Range("P7").Value = 0.02733363438148323
$ws.Range("S7").Value = 0.002462935221472211
$ws.Range("T7").Value = 0.002462935221472212
$ws.Range("I8").Value = 0.09010639372350319
$ws.Range("J8").Value = 0.09010639372350321
$ws.Range("M8").Value = 7.880893333333333
$ws.Range("N8").Value = 23.64268
$ws.Range("O8").Value = 0.281520346184098
$ws.Range("P8").Value = 0.281520346184098
$ws.Range("Q8").Value = 19.47535817605333
$ws.Range("R8").Value = 175.27822358448
$ws.Range("S8").Value = 0.02536678315444126
$ws.Range("T8").Value = 0.02536678315444126
$ws.Range("I9").Value = 0.09010639372350319
$ws.Range("J9").Value = 0.09010639372350321
$ws.Range("M9").Value = 0.183999
$ws.Range("N9").Value = 0.551997
$ws.Range("O9").Value = 0.006572790670625477
$ws.Range("P9").Value = 0.006572790670625476
$ws.Range("Q9").Value = 0.454700536788
$ws.Range("R9").Value = 4.092304831092
$ws.Range("S9").Value = 0.0005922504640295479
$ws.Range("T9").Value = 0.0005922504640295479
$ws.Range("G10").Value = 1.148663
$ws.Range("H10").Value = 3.445989
$ws.Range("I10").Value = 0.04188304383987305
$ws.Range("J10").Value = 0.04188304383987305
$ws.Range("M10").Value = 19.163974
$ws.Range("N10").Value = 57.491922
$ws.Range("O10").Value = 0.6845732287637933
$ws.Range("P10").Value = 0.6845732287637933
$ws.Range("Q10").Value = 22.012947866762
$ws.Range("R10").Value = 198.116530800858
$ws.Range("S10").Value = 0.0286720105519174
$ws.Range("T10").Value = 0.0286720105519174
$ws.Range("G11").Value = 1.148663
$ws.Range("H11").Value = 3.445989
$ws.Range("I11").Value = 0.04188304383987305
$ws.Range("J11").Value = 0.04188304383987305
$ws.Range("O11").Value = 0.02733363438148322
$ws.Range("P11").Value = 0.02733363438148323
$ws.Range("Q11").Value = 0.8789328056769999
$ws.Range("R11").Value = 7.910395251093
$ws.Range("S11").Value = 0.001144815807102723
$ws.Range("T11").Value = 0.001144815807102723
$ws.Range("G12").Value = 1.148663
$ws.Range("H12").Value = 3.445989
$ws.Range("I12").Value = 0.04188304383987305
$ws.Range("J12").Value = 0.04188304383987305
$ws.Range("M12").Value = 7.880893333333333
$ws.Range("N12").Value = 23.64268
$ws.Range("O12").Value = 0.281520346184098
$ws.Range("P12").Value = 0.281520346184098
$ws.Range("Q12").Value = 9.052490578946665
$ws.Range("R12").Value = 81.47241521052
$ws.Range("S12").Value = 0.01179092900104481
$ws.Range("T12").Value = 0.01179092900104481
$ws.Range("G13").Value = 1.148663
$ws.Range("H13").Value = 3.445989
$ws.Range("I13").Value = 0.04188304383987305
$ws.Range("J13").Value = 0.04188304383987305
$ws.Range("M13").Value = 0.183999
$ws.Range("N13").Value = 0.551997
$ws.Range("O13").Value = 0.006572790670625477
$ws.Range("P13").Value = 0.006572790670625476
$ws.Range("Q13").Value = 0.211352843337
$ws.Range("R13").Value = 1.902175590033
$ws.Range("S13").Value = 0.0002752884798081154
$ws.Range("T13").Value = 0.0002752884798081154
$ws.Range("G14").Value = 0.07412966666666666
$ws.Range("H14").Value = 0.222389
$ws.Range("I14").Value = 0.002702947756509242
$ws.Range("J14").Value = 0.002702947756509243
$ws.Range("M14").Value = 19.163974
$ws.Range("N14").Value = 57.491922
$ws.Range("O14").Value = 0.6845732287637933
$ws.Range("P14").Value = 0.6845732287637933
$ws.Range("Q14").Value = 1.420619004628667
$ws.Range("R14").Value = 12.785571041658
$ws.Range("S14").Value = 0.001850365672853384
$ws.Range("T14").Value = 0.001850365672853384
$ws.Range("G15").Value = 0.07412966666666666
$ws.Range("H15").Value = 0.222389
$ws.Range("I15").Value = 0.002702947756509242
$ws.Range("J15").Value = 0.002702947756509243
$ws.Range("O15").Value = 0.02733363438148322
$ws.Range("P15").Value = 0.02733363438148323
$ws.Range("Q15").Value = 0.05672246421033333
$ws.Range("R15").Value = 0.510502177893
$ws.Range("S15").Value = 0.00007388138572867397
$ws.Range("T15").Value = 0.00007388138572867398
$ws.Range("G16").Value = 0.07412966666666666
$ws.Range("H16").Value = 0.222389
$ws.Range("I16").Value = 0.002702947756509242
$ws.Range("J16").Value = 0.002702947756509243
$ws.Range("M16").Value = 7.880893333333333
$ws.Range("N16").Value = 23.64268
$ws.Range("O16").Value = 0.281520346184098
$ws.Range("P16").Value = 0.281520346184098
$ws.Range("Q16").Value = 0.5842079958355555
$ws.Range("R16").Value = 5.257871962519999
$ws.Range("S16").Value = 0.0007609347881300129
$ws.Range("T16").Value = 0.000760934788130013
$ws.Range("G17").Value = 0.07412966666666666
$ws.Range("H17").Value = 0.222389
$ws.Range("I17").Value = 0.002702947756509242
$ws.Range("J17").Value = 0.002702947756509243
$ws.Range("M17").Value = 0.183999
$ws.Range("N17").Value = 0.551997
$ws.Range("O17").Value = 0.006572790670625477
$ws.Range("P17").Value = 0.006572790670625476
$ws.Range("Q17").Value = 0.013639784537
$ws.Range("R17").Value = 0.122758060833
$ws.Range("S17").Value = 0.00001776590979717201
$ws.Range("T17").Value = 0.00001776590979717201
